$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell updates from the crypto price refresh.
# Column D holds numeric-looking price text (must stay text, incl. trailing
# zeros / multi-dot thousands separators), so we force NumberFormat="@" there.
# Columns B, C, E are unambiguous text already (names/URLs, or padded % strings).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.942.19"
$ws.Range("E2").Value = "  -0.07%  "
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.354.31"
$ws.Range("E3").Value = "  -0.12%  "
# Row 4
$ws.Range("E4").Value = "  +0.09%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.53"
$ws.Range("E5").Value = "  +0.23%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.668"
$ws.Range("E6").Value = "  -2.72%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.04"
$ws.Range("E7").Value = "  -2.89%  "
# Row 8
$ws.Range("E8").Value = "  +0.01%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.599"
$ws.Range("E9").Value = "  -3.40%  "
# Row 10
$ws.Range("E10").Value = "  +0.04%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.87"
$ws.Range("E11").Value = "  +4.75%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.21"
$ws.Range("E12").Value = "  -0.08%  "
# Row 13
$ws.Range("E13").Value = "  +0.69%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.22"
$ws.Range("E14").Value = "  -2.72%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.705.29"
$ws.Range("E15").Value = "  -0.07%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.22"
$ws.Range("E16").Value = "  -2.34%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.906"
$ws.Range("E17").Value = "  -0.72%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.357.17"
$ws.Range("E18").Value = "  +0.15%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.869.83"
# Row 20
$ws.Range("E20").Value = "  +0.94%  "
# Row 21
$ws.Range("E21").Value = "  +0.30%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "78.13"
$ws.Range("E22").Value = "  +0.87%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.94"
$ws.Range("E23").Value = "  -2.26%  "
# Row 24
$ws.Range("E24").Value = "  +0.15%  "
# Row 25
$ws.Range("E25").Value = "  +2.39%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.84"
$ws.Range("E26").Value = "  -0.84%  "
# Row 27
$ws.Range("E27").Value = "  -0.19%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.48"
$ws.Range("E28").Value = "  -2.80%  "
# Row 29
$ws.Range("B29").Value = "Monero"
$ws.Range("C29").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "176.53"
$ws.Range("E29").Value = "  +0.91%  "
# Row 30
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.20"
$ws.Range("E30").Value = "  -1.80%  "
# Row 32
$ws.Range("E32").Value = "  +0.38%  "
# Row 33
$ws.Range("E33").Value = "  -1.27%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0750"
$ws.Range("E34").Value = "  -0.82%  "
# Row 35
$ws.Range("B35").Value = "InternetComputer(DFINITY)"
$ws.Range("C35").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.40"
$ws.Range("E35").Value = "  +0.31%  "
# Row 36
$ws.Range("B36").Value = "Filecoin"
$ws.Range("C36").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.09"
$ws.Range("E36").Value = "  -4.41%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.83"
$ws.Range("E37").Value = "  +2.60%  "
# Row 38
$ws.Range("E38").Value = "  +0.42%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.43"
$ws.Range("E39").Value = "  +0.79%  "
# Row 40
$ws.Range("E40").Value = "  -3.27%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.47"
$ws.Range("E41").Value = "  +14.43%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "64.78"
$ws.Range("E42").Value = "  +13.67%  "
# Row 43
$ws.Range("E43").Value = "  -0.60%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.20"
$ws.Range("E44").Value = "  -0.86%  "
# Row 45
$ws.Range("E45").Value = "  -2.73%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.201"
$ws.Range("E46").Value = "  -6.08%  "
# Row 47
$ws.Range("E47").Value = "  -0.06%  "
# Row 48
$ws.Range("E48").Value = "  -2.07%  "
# Row 49
$ws.Range("E49").Value = "  -6.19%  "
# Row 50
$ws.Range("E50").Value = "  -1.95%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "98.64"
$ws.Range("E51").Value = "  -1.81%  "
